$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate new rows 13-24 with cell values ---

# Row 13: RegisterWithUnmatchedPasswords
$ws.Range("A13").Value = 'RegisterWithUnmatchedPasswords'
$ws.Range("B13").Value = 'Iliya'
$ws.Range("C13").Value = 'Iliev'
$ws.Range("D13").Value = 'true true true'
$ws.Range("E13").Value = 'false false true'
$ws.Range("F13").Value = 'Bulgaria'
$ws.Range("G13").Value = '3'
$ws.Range("H13").Value = '23'
$ws.Range("I13").Value = '1987'
$ws.Range("J13").Value = '0897675645'
$ws.Range("K13").Value = 'lichkata456'
$ws.Range("L13").Value = 'lichkata456@abv.bg'
$ws.Range("M13").Value = 'C:\Users\Iliya\Desktop\photo.jpeg'
$ws.Range("N13").Value = 'ALA BALA'
$ws.Range("O13").Value = '12345678'
$ws.Range("P13").Value = '12345679'

# Row 14: RegisterWithoutLastNameAndPhone
$ws.Range("A14").Value = 'RegisterWithoutLastNameAndPhone'
$ws.Range("B14").Value = 'Iliya'
$ws.Range("C14").Value = 'String.Empty'
$ws.Range("D14").Value = 'true true true'
$ws.Range("E14").Value = 'false false true'
$ws.Range("F14").Value = 'Bulgaria'
$ws.Range("G14").Value = '3'
$ws.Range("H14").Value = '23'
$ws.Range("I14").Value = '1987'
$ws.Range("J14").Value = 'String.Empty'
$ws.Range("K14").Value = 'lichkata456'
$ws.Range("L14").Value = 'lichkata456@abv.bg'
$ws.Range("M14").Value = 'C:\Users\Iliya\Desktop\photo.jpeg'
$ws.Range("N14").Value = 'ALA BALA'
$ws.Range("O14").Value = '12345678'
$ws.Range("P14").Value = '12345678'

# Row 15: RegisterWithoutLastNameAndUsername
$ws.Range("A15").Value = 'RegisterWithoutLastNameAndUsername'
$ws.Range("B15").Value = 'Iliya'
$ws.Range("C15").Value = 'String.Empty'
$ws.Range("D15").Value = 'true true true'
$ws.Range("E15").Value = 'false false true'
$ws.Range("F15").Value = 'Bulgaria'
$ws.Range("G15").Value = '3'
$ws.Range("H15").Value = '23'
$ws.Range("I15").Value = '1987'
$ws.Range("J15").Value = '0897675645'
$ws.Range("K15").Value = 'String.Empty'
$ws.Range("L15").Value = 'lichkata456@abv.bg'
$ws.Range("M15").Value = 'C:\Users\Iliya\Desktop\photo.jpeg'
$ws.Range("N15").Value = 'ALA BALA'
$ws.Range("O15").Value = '12345678'
$ws.Range("P15").Value = '12345678'

# Row 16: RegisterWithoutLastNameAndEmail
$ws.Range("A16").Value = 'RegisterWithoutLastNameAndEmail'
$ws.Range("B16").Value = 'Iliya'
$ws.Range("C16").Value = 'String.Empty'
$ws.Range("D16").Value = 'true true true'
$ws.Range("E16").Value = 'false false true'
$ws.Range("F16").Value = 'Bulgaria'
$ws.Range("G16").Value = '3'
$ws.Range("H16").Value = '23'
$ws.Range("I16").Value = '1987'
$ws.Range("J16").Value = '0897675645'
$ws.Range("K16").Value = 'lichkata456'
$ws.Range("L16").Value = 'String.Empty'
$ws.Range("M16").Value = 'C:\Users\Iliya\Desktop\photo.jpeg'
$ws.Range("N16").Value = 'ALA BALA'
$ws.Range("O16").Value = '12345678'
$ws.Range("P16").Value = '12345678'

# Row 17: RegisterWithoutLastNameAndWrongEmail
$ws.Range("A17").Value = 'RegisterWithoutLastNameAndWrongEmail'
$ws.Range("B17").Value = 'Iliya'
$ws.Range("C17").Value = 'String.Empty'
$ws.Range("D17").Value = 'true true true'
$ws.Range("E17").Value = 'false false true'
$ws.Range("F17").Value = 'Bulgaria'
$ws.Range("G17").Value = '3'
$ws.Range("H17").Value = '23'
$ws.Range("I17").Value = '1987'
$ws.Range("J17").Value = '0897675645'
$ws.Range("K17").Value = 'lichkata456'
$ws.Range("L17").Value = 'lichkata456@abv.'
$ws.Range("M17").Value = 'C:\Users\Iliya\Desktop\photo.jpeg'
$ws.Range("N17").Value = 'ALA BALA'
$ws.Range("O17").Value = '12345678'
$ws.Range("P17").Value = '12345678'

# Row 18: RegisterWithoutHobbiesAndUsername
$ws.Range("A18").Value = 'RegisterWithoutHobbiesAndUsername'
$ws.Range("B18").Value = 'Iliya'
$ws.Range("C18").Value = 'Iliev'
$ws.Range("D18").Value = 'true true true'
$ws.Range("E18").Value = 'false false false'
$ws.Range("F18").Value = 'Bulgaria'
$ws.Range("G18").Value = '3'
$ws.Range("H18").Value = '23'
$ws.Range("I18").Value = '1987'
$ws.Range("J18").Value = '0897675645'
$ws.Range("K18").Value = 'String.Empty'
$ws.Range("L18").Value = 'lichkata456@abv.bg'
$ws.Range("M18").Value = 'C:\Users\Iliya\Desktop\photo.jpeg'
$ws.Range("N18").Value = 'ALA BALA'
$ws.Range("O18").Value = '12345678'
$ws.Range("P18").Value = '12345678'

# Row 19: RegisterWithoutUsernameAndPassword
$ws.Range("A19").Value = 'RegisterWithoutUsernameAndPassword'
$ws.Range("B19").Value = 'Iliya'
$ws.Range("C19").Value = 'Iliev'
$ws.Range("D19").Value = 'true true true'
$ws.Range("E19").Value = 'false true false'
$ws.Range("F19").Value = 'Bulgaria'
$ws.Range("G19").Value = '3'
$ws.Range("H19").Value = '23'
$ws.Range("I19").Value = '1987'
$ws.Range("J19").Value = '0897675645'
$ws.Range("K19").Value = 'String.Empty'
$ws.Range("L19").Value = 'lichkata456@abv.bg'
$ws.Range("M19").Value = 'C:\Users\Iliya\Desktop\photo.jpeg'
$ws.Range("N19").Value = 'String.Empty'
$ws.Range("O19").Value = '12345678'
$ws.Range("P19").Value = '12345678'

# Row 20: RegisterWithNegativePhone
$ws.Range("A20").Value = 'RegisterWithNegativePhone'
$ws.Range("B20").Value = 'Iliya'
$ws.Range("C20").Value = 'Iliev'
$ws.Range("D20").Value = 'true true true'
$ws.Range("E20").Value = 'false false true'
$ws.Range("F20").Value = 'Bulgaria'
$ws.Range("G20").Value = '3'
$ws.Range("H20").Value = '23'
$ws.Range("I20").Value = '1987'
$ws.Range("J20").Value = '-0897675646'
$ws.Range("K20").Value = 'lichkata457'
$ws.Range("L20").Value = 'lichkata456@abv.bg'
$ws.Range("M20").Value = 'C:\Users\Iliya\Desktop\photo.jpeg'
$ws.Range("N20").Value = 'ALA BALA'
$ws.Range("O20").Value = '12345678'
$ws.Range("P20").Value = '12345678'

# Row 21: RegisterWithNegativePhoneAndWrongEmail
$ws.Range("A21").Value = 'RegisterWithNegativePhoneAndWrongEmail'
$ws.Range("B21").Value = 'Iliya'
$ws.Range("C21").Value = 'Iliev'
$ws.Range("D21").Value = 'true true true'
$ws.Range("E21").Value = 'false false true'
$ws.Range("F21").Value = 'Bulgaria'
$ws.Range("G21").Value = '3'
$ws.Range("H21").Value = '23'
$ws.Range("I21").Value = '1987'
$ws.Range("J21").Value = '-0897675646'
$ws.Range("K21").Value = 'lichkata457'
$ws.Range("L21").Value = 'lichkata456@abv.'
$ws.Range("M21").Value = 'C:\Users\Iliya\Desktop\photo.jpeg'
$ws.Range("N21").Value = 'ALA BALA'
$ws.Range("O21").Value = '12345678'
$ws.Range("P21").Value = '12345678'

# Row 22: RegisterWithWrongEmailAndMismatchedPasswords
$ws.Range("A22").Value = 'RegisterWithWrongEmailAndMismatchedPasswords'
$ws.Range("B22").Value = 'Iliya'
$ws.Range("C22").Value = 'Iliev'
$ws.Range("D22").Value = 'true true true'
$ws.Range("E22").Value = 'false false true'
$ws.Range("F22").Value = 'Bulgaria'
$ws.Range("G22").Value = '3'
$ws.Range("H22").Value = '23'
$ws.Range("I22").Value = '1987'
$ws.Range("J22").Value = '0897675646'
$ws.Range("K22").Value = 'lichkata457'
$ws.Range("L22").Value = 'lichkata456@abv.'
$ws.Range("M22").Value = 'C:\Users\Iliya\Desktop\photo.jpeg'
$ws.Range("N22").Value = 'ALA BALA'
$ws.Range("O22").Value = '12345678'
$ws.Range("P22").Value = '123456789'

# Row 23: RegisterWithWrongEmailAndMismatchedAndShortPasswords
$ws.Range("A23").Value = 'RegisterWithWrongEmailAndMismatchedAndShortPasswords'
$ws.Range("B23").Value = 'Iliya'
$ws.Range("C23").Value = 'Iliev'
$ws.Range("D23").Value = 'true true true'
$ws.Range("E23").Value = 'false false true'
$ws.Range("F23").Value = 'Bulgaria'
$ws.Range("G23").Value = '3'
$ws.Range("H23").Value = '23'
$ws.Range("I23").Value = '1987'
$ws.Range("J23").Value = '0897675646'
$ws.Range("K23").Value = 'lichkata457'
$ws.Range("L23").Value = 'lichkata456@abv.'
$ws.Range("M23").Value = 'C:\Users\Iliya\Desktop\photo.jpeg'
$ws.Range("N23").Value = 'ALA BALA'
$ws.Range("O23").Value = '1234567'
$ws.Range("P23").Value = '123456789'

# Row 24: RegisterAlreadyRegisteredUser
$ws.Range("A24").Value = 'RegisterAlreadyRegisteredUser'
$ws.Range("B24").Value = 'Iliya'
$ws.Range("C24").Value = 'Iliev'
$ws.Range("D24").Value = 'true true true'
$ws.Range("E24").Value = 'true true true'
$ws.Range("F24").Value = 'Bulgaria'
$ws.Range("G24").Value = '3'
$ws.Range("H24").Value = '23'
$ws.Range("I24").Value = '1987'
$ws.Range("J24").Value = '0897675645'
$ws.Range("K24").Value = 'lichkata456'
$ws.Range("L24").Value = 'lichkata456@abv.bg'
$ws.Range("M24").Value = 'C:\Users\Iliya\Desktop\photo.jpeg'
$ws.Range("N24").Value = 'ALA BALA'
$ws.Range("O24").Value = '12345678'
$ws.Range("P24").Value = '12345678'

# --- Apply Hyperlink style + hyperlink relationships on column L for rows with an email ---

$ws.Hyperlinks.Add($ws.Range("L13"), 'mailto:lichkata456@abv.bg') | Out-Null
$ws.Range("L13").Style = "Hyperlink"
$ws.Range("L13").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("L14"), 'mailto:lichkata456@abv.bg') | Out-Null
$ws.Range("L14").Style = "Hyperlink"
$ws.Range("L14").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("L15"), 'mailto:lichkata456@abv.bg') | Out-Null
$ws.Range("L15").Style = "Hyperlink"
$ws.Range("L15").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("L17"), 'mailto:lichkata456@abv.') | Out-Null
$ws.Range("L17").Style = "Hyperlink"
$ws.Range("L17").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("L18"), 'mailto:lichkata456@abv.bg') | Out-Null
$ws.Range("L18").Style = "Hyperlink"
$ws.Range("L18").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("L19"), 'mailto:lichkata456@abv.bg') | Out-Null
$ws.Range("L19").Style = "Hyperlink"
$ws.Range("L19").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("L20"), 'mailto:lichkata456@abv.bg') | Out-Null
$ws.Range("L20").Style = "Hyperlink"
$ws.Range("L20").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("L21"), 'mailto:lichkata456@abv.') | Out-Null
$ws.Range("L21").Style = "Hyperlink"
$ws.Range("L21").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("L22"), 'mailto:lichkata456@abv.') | Out-Null
$ws.Range("L22").Style = "Hyperlink"
$ws.Range("L22").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("L23"), 'mailto:lichkata456@abv.') | Out-Null
$ws.Range("L23").Style = "Hyperlink"
$ws.Range("L23").NumberFormat = "@"

$ws.Hyperlinks.Add($ws.Range("L24"), 'mailto:lichkata456@abv.bg') | Out-Null
$ws.Range("L24").Style = "Hyperlink"
$ws.Range("L24").NumberFormat = "@"

# --- Column A width adjustment (manual resize, bestFit removed) ---
$ws.Columns("A").ColumnWidth = 34.42

# --- Selection / view state ---
$ws.Range("C27").Select() | Out-Null

